$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Preserve text (inline-string) cell type for numeric-looking values by
# forcing Text number format before assignment (otherwise COM coerces
# numeric-looking strings into real numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.19'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '25.53'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.115'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05582'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.471'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.016'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8184'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8406'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0005950'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '9OneONE'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1334'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02876'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09371'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001509'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '13BitForexTokenBF'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.006218'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '14TigerCashTCH'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.523'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '15LEOLEO'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.022'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '16BTSETokenBTSE'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.3179'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '17BitpandaEcosystemTokenBEST'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'MandalaExchangeToken'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06952'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '18MandalaExchangeTokenMDX'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03225'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.744'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04707'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001250'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004610'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009700'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03656'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '40KickTokenKICKBestin24h'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1051'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002506'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007643'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005313'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1335'
